# Applies the "duplicate columns A:D into E:H" edit to every worksheet in the
# workbook (Orange, Alstom, ... one block per existing sheet), and extends the
# header merge (A1:D1) with a matching merge (E1:H1).
#
# Rows 2-10 contain 4 populated columns (label / % change / high / low) that
# get mirrored into E:H.
# Rows 11-14 (MM20 / MM50 / MM100 / RSI14) only have the label (col A) and the
# value (col D) populated - columns B/C (and their mirrors F/G) stay blank.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Row 1: header date - mirror A1 into E1 and merge E1:H1 to match A1:D1.
    $ws.Cells.Item(1, 5).Value2 = $ws.Cells.Item(1, 1).Value2
    $ws.Range("E1:H1").Merge()

    # Rows 2-10: full rows -> copy columns A,B,C,D into E,F,G,H.
    for ($r = 2; $r -le 10; $r++) {
        for ($c = 1; $c -le 4; $c++) {
            $srcValue = $ws.Cells.Item($r, $c).Value2
            $dest = $ws.Cells.Item($r, $c + 4)
            if ($c -eq 2) {
                # Column B holds a literal "+x.xx%"/"-x.xx%" text string.
                # Writing it straight to .Value2 makes Excel auto-convert it
                # to a numeric percentage, so briefly force a text format
                # while assigning, then clear the format again so the
                # destination cell is left with the default (unstyled)
                # formatting, same as the rest of that column.
                $dest.NumberFormat = "@"
                $dest.Value2 = $srcValue
                $dest.ClearFormats()
            } else {
                $dest.Value2 = $srcValue
            }
        }
    }

    # Rows 11-14: only label (A) and value (D) are populated; B/C (and the
    # mirrored F/G) remain blank.
    for ($r = 11; $r -le 14; $r++) {
        $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 1).Value2
        $ws.Cells.Item($r, 6).Value2 = $null
        $ws.Cells.Item($r, 7).Value2 = $null
        $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 4).Value2

        $ws.Cells.Item($r, 2).Value2 = $null
        $ws.Cells.Item($r, 3).Value2 = $null
    }
}
